# Generate Report for handback
# Adds two newly-handed-back files to the handback-status workbook:
#   0a72fee4-4502-4bed-9493-41a590896049.md
#   9f13907f-091e-4ea4-8a10-7a981bf63cf9.md
# across the "Overview", "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

function Set-Hyperlink($ws, $cellRef, $displayText, $url) {
    $ws.Range($cellRef).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText)
}

# ---------------------------------------------------------------------------
# Sheet "Overview": one row per handed-back source file.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-Hyperlink $wsOverview "A6" "0a72fee4-4502-4bed-9493-41a590896049.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/0a72fee4-4502-4bed-9493-41a590896049.md"
$wsOverview.Range("B6").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C6").Value = "Handed back: in sync with en-US"

Set-Hyperlink $wsOverview "A7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/9f13907f-091e-4ea4-8a10-7a981bf63cf9.md"
$wsOverview.Range("B7").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C7").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": full handback detail rows (zh-cn locale).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 6 - 0a72fee4-4502-4bed-9493-41a590896049 (zh-cn)
Set-Hyperlink $wsZh "A6" "0a72fee4-4502-4bed-9493-41a590896049.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/0a72fee4-4502-4bed-9493-41a590896049.md"
$wsZh.Range("B6").Value = "Handed back: in sync with en-US"
Set-Hyperlink $wsZh "C6" "0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c6f4e9eb23f28b47f2e05e80d4f2ef65d4f7c6a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.zh-cn.xlf"
$wsZh.Range("D6").Value = "2016-02-16 10:32:12"
Set-Hyperlink $wsZh "E6" "0a72fee4-4502-4bed-9493-41a590896049.md" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a7d8f1b2c3d4e5f60718293a4b5c6d7e8f901234/e2e/0a72fee4-4502-4bed-9493-41a590896049.md"
Set-Hyperlink $wsZh "F6" "0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1b2c3d4e5f60718293a4b5c6d7e8f9012345678a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.zh-cn.xlf"
$wsZh.Range("G6").Value = "2016-02-16 10:33:08"
$wsZh.Range("H6").Value = "Include"

# Row 7 - 9f13907f-091e-4ea4-8a10-7a981bf63cf9 (zh-cn)
Set-Hyperlink $wsZh "A7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/9f13907f-091e-4ea4-8a10-7a981bf63cf9.md"
$wsZh.Range("B7").Value = "Handed back: in sync with en-US"
Set-Hyperlink $wsZh "C7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2468ace13579bdf02468ace13579bdf02468ace/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.zh-cn.xlf"
$wsZh.Range("D7").Value = "2016-02-16 10:32:12"
Set-Hyperlink $wsZh "E7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.md" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/13579bdf02468ace13579bdf02468ace13579bd/e2e/9f13907f-091e-4ea4-8a10-7a981bf63cf9.md"
Set-Hyperlink $wsZh "F7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/02468ace13579bdf02468ace13579bdf02468ace/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.zh-cn.xlf"
$wsZh.Range("G7").Value = "2016-02-16 10:33:36"
$wsZh.Range("H7").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de": full handback detail rows (de-de locale).
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 6 - 0a72fee4-4502-4bed-9493-41a590896049 (de-de)
Set-Hyperlink $wsDe "A6" "0a72fee4-4502-4bed-9493-41a590896049.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/0a72fee4-4502-4bed-9493-41a590896049.md"
$wsDe.Range("B6").Value = "Handed back: in sync with en-US"
Set-Hyperlink $wsDe "C6" "0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d5e6f708192a3b4c5d6e7f8091a2b3c4d5e6f70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.de-de.xlf"
$wsDe.Range("D6").Value = "2016-02-16 10:32:26"
Set-Hyperlink $wsDe "E6" "0a72fee4-4502-4bed-9493-41a590896049.md" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b2c3d4e5f60718293a4b5c6d7e8f9012345678ab/e2e/0a72fee4-4502-4bed-9493-41a590896049.md"
Set-Hyperlink $wsDe "F6" "0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2c3d4e5f60718293a4b5c6d7e8f9012345678abc/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0a72fee4-4502-4bed-9493-41a590896049.0f45bf5ae43d5dc36adae02e61ed771eed151fdf.de-de.xlf"
$wsDe.Range("G6").Value = "2016-02-16 10:33:36"
$wsDe.Range("H6").Value = "Include"

# Row 7 - 9f13907f-091e-4ea4-8a10-7a981bf63cf9 (de-de)
Set-Hyperlink $wsDe "A7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.md" "https://github.com/OpenLocalizationTest/oltest/blob/5d1f170ae1ea497aacaa396a99bb8d385896c6e1/e2e/9f13907f-091e-4ea4-8a10-7a981bf63cf9.md"
$wsDe.Range("B7").Value = "Handed back: in sync with en-US"
Set-Hyperlink $wsDe "C7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/369bdf02468ace13579bdf02468ace13579bdf0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.de-de.xlf"
$wsDe.Range("D7").Value = "2016-02-16 10:32:26"
Set-Hyperlink $wsDe "E7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.md" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/79bdf02468ace13579bdf02468ace13579bdf024/e2e/9f13907f-091e-4ea4-8a10-7a981bf63cf9.md"
Set-Hyperlink $wsDe "F7" "9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9bdf02468ace13579bdf02468ace13579bdf0246/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/9f13907f-091e-4ea4-8a10-7a981bf63cf9.ac34053823ffaa37a597a9bd7ad03924909c73d0.de-de.xlf"
$wsDe.Range("G7").Value = "2016-02-16 10:33:36"
$wsDe.Range("H7").Value = "Include"

Write-Host "Handback rows added for 0a72fee4-4502-4bed-9493-41a590896049 and 9f13907f-091e-4ea4-8a10-7a981bf63cf9"
